$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge "<TableView> <and >[bookmark]<has a >" runs into a
# single run "<TableView>< and has a >", dropping the _GoBack bookmark.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(" and has a ", $true, $false, $false, $false, $false, $true, 1, $false, " and has a ", 2)

# ---------------------------------------------------------------------
# Change 2: restructure the "Navigation:" section --
#   * delete the blank paragraph right before "Navigation:"
#   * rename that heading paragraph "Navigation:" -> "Core Data:"
#   * insert a new body paragraph explaining Core Data usage
#   * insert a new "Navigation:" heading paragraph after it
#   * move the _GoBack bookmark into the paragraph that used to follow
#     "Navigation:"
# ---------------------------------------------------------------------

# Locate "Navigation:" and remove the empty paragraph preceding it.
$rng = $d.Content
$null = $rng.Find.Execute("Navigation:")
$navPara = $rng.Paragraphs(1)
$blankBefore = $navPara.Previous()
$blankBefore.Range.Delete()

# Re-locate "Navigation:" (paragraph objects go stale after the edit)
# and rename it to "Core Data:".
$rng = $d.Content
$null = $rng.Find.Execute("Navigation:")
$headingPara = $rng.Paragraphs(1)
$headingPara.Range.Text = "Core Data:"

# Re-locate "Core Data:" and append two placeholder paragraphs after it
# (their inherited formatting will be overwritten below).
$rng = $d.Content
$null = $rng.Find.Execute("Core Data:")
$headingPara = $rng.Paragraphs(1)
$headingPara.Range.InsertParagraphAfter()

$rng = $d.Content
$null = $rng.Find.Execute("Core Data:")
$headingPara = $rng.Paragraphs(1)
$bodyPlaceholder = $headingPara.Next()
$bodyPlaceholder.Range.InsertParagraphAfter()

$rng = $d.Content
$null = $rng.Find.Execute("Core Data:")
$headingPara = $rng.Paragraphs(1)
$bodyPlaceholder = $headingPara.Next()
$navPlaceholder = $bodyPlaceholder.Next()

$bodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:right="-1800"/><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:tab/><w:t xml:space="preserve">We used the Core Data persistence model, and it only includes one </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t>entity which</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve"> is Trail. We fetch the data in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t>TableViewControllers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t>, and fro</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve">m </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t>there</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve"> pass the data to other </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve">View Controllers instead of </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t>fetching it</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/></w:rPr><w:t xml:space="preserve"> in every View Controller.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$navXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:right="-1800"/><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS" w:cs="Arial Unicode MS"/><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t>Navigation:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$bodyPlaceholder.Range.InsertXML($bodyXml)
$navPlaceholder.Range.InsertXML($navXml)

# Re-locate the new "Navigation:" heading and add the _GoBack bookmark
# to the paragraph that now follows it.
$rng = $d.Content
$null = $rng.Find.Execute("Navigation:")
$navPara = $rng.Paragraphs(1)
$afterNav = $navPara.Next()
$null = $d.Bookmarks.Add("_GoBack", $afterNav.Range)

Write-Output "done"
